# Files upload sections updated
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A header/meter-id -> user-id rename
$ws.Range("A1").Value = "userId"

# Sample ids simplified
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 2

# Reset the view: scroll back to column A and move the active selection to A4
$ws.Range("A4").Select()
